$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the three "check/validate" test-step names to the new
# machine-style step identifiers (B6:B8 on the "Test Flow" sheet).
$ws.Range("B6").Value = "CHECK_INTERLOCK"
$ws.Range("B7").Value = "HIPOT_TEST"
$ws.Range("B8").Value = "HIPOT_RESET"

# Row 6 ("CHECK_INTERLOCK") is no longer enabled by default.
$ws.Range("E6").Value = $false

# Reflect the user's on-screen selection after the edit: the A6:E8
# block (the three test-step rows) selected, anchored at E8.
$ws.Range("A6:E8").Select()
